$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4628.0625
$ws.Range("I74").Value = 4912.4165
$ws.Range("J74").Value = 3775
$ws.Range("K74").Value = 4912.4165
$ws.Range("L74").Value = 3775
$ws.Range("M74").Value = -3976.4165
$ws.Range("N74").Value = -5647

$ws.Range("H77").Value = 4628.0625
$ws.Range("I77").Value = 4912.4165
$ws.Range("J77").Value = 3775
$ws.Range("K77").Value = 24562.0825
$ws.Range("L77").Value = 18875
$ws.Range("M77").Value = -19882.0825
$ws.Range("N77").Value = -28235

$ws.Range("H125").Value = 795.6667
$ws.Range("J125").Value = 840
$ws.Range("L125").Value = 7560
$ws.Range("N125").Value = -12480

$ws.Range("H132").Value = 16033752
$ws.Range("I132").Value = 19058570
$ws.Range("J132").Value = 2218.4
$ws.Range("K132").Value = 57175710
$ws.Range("L132").Value = 6655.200000000001
$ws.Range("M132").Value = -57173180
$ws.Range("N132").Value = -11715.2

$ws.Range("H137").Value = 1109.6293
$ws.Range("I137").Value = 948.3570999999999
$ws.Range("J137").Value = 1383.303
$ws.Range("K137").Value = 2845.0713
$ws.Range("L137").Value = 4149.909000000001
$ws.Range("M137").Value = -295.0712999999996
$ws.Range("N137").Value = -9249.909

$ws.Range("H138").Value = 12822.965
$ws.Range("I138").Value = 3914.3584
$ws.Range("J138").Value = 27577.844
$ws.Range("K138").Value = 11743.0752
$ws.Range("L138").Value = 82733.53200000001
$ws.Range("M138").Value = -6603.075199999999
$ws.Range("N138").Value = -93013.53200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 40712
$ws.Range("J7").Value = 40712
$ws.Range("L7").Value = 40712
$ws.Range("N7").Value = -40940

$ws.Range("H61").Value = 4104.674
$ws.Range("I61").Value = 4532.4414
$ws.Range("J61").Value = 2892.6667
$ws.Range("K61").Value = 4532.4414
$ws.Range("L61").Value = 2892.6667
$ws.Range("M61").Value = -4320.4414
$ws.Range("N61").Value = -3316.6667

$ws.Range("H64").Value = 17500
$ws.Range("I64").Value = 5000
$ws.Range("J64").Value = 30000
$ws.Range("K64").Value = 5000
$ws.Range("L64").Value = 30000
$ws.Range("M64").Value = -4752
$ws.Range("N64").Value = -30496

$ws.Range("H67").Value = 17500
$ws.Range("I67").Value = 5000
$ws.Range("J67").Value = 30000
$ws.Range("K67").Value = 5000
$ws.Range("L67").Value = 30000
$ws.Range("M67").Value = -4142
$ws.Range("N67").Value = -31716

$ws.Range("H109").Value = 20000
$ws.Range("J109").Value = 20000
$ws.Range("L109").Value = 20000
$ws.Range("N109").Value = -22774

$ws.Range("H122").Value = 34483416
$ws.Range("I122").Value = 37037696
$ws.Range("J122").Value = 660
$ws.Range("K122").Value = 111113088
$ws.Range("L122").Value = 1980
$ws.Range("M122").Value = -111110638
$ws.Range("N122").Value = -6880

$ws.Range("H136").Value = 4104.674
$ws.Range("I136").Value = 4532.4414
$ws.Range("J136").Value = 2892.6667
$ws.Range("K136").Value = 13597.3242
$ws.Range("L136").Value = 8678.000100000001
$ws.Range("M136").Value = -11047.3242
$ws.Range("N136").Value = -13778.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 22685.715
$ws.Range("J2").Value = 22685.715
$ws.Range("L2").Value = 22685.715
$ws.Range("N2").Value = -22911.715

$ws.Range("H62").Value = 27417.715
$ws.Range("J62").Value = 28987.334
$ws.Range("L62").Value = 28987.334
$ws.Range("N62").Value = -30359.334

$ws.Range("H65").Value = 27417.715
$ws.Range("J65").Value = 28987.334
$ws.Range("L65").Value = 86962.00199999999
$ws.Range("N65").Value = -93826.00199999999

$ws.Range("H108").Value = 27263

$ws.Range("H134").Value = 16693623
$ws.Range("I134").Value = 20866128
$ws.Range("J134").Value = 3600
$ws.Range("K134").Value = 62598384
$ws.Range("L134").Value = 10800
$ws.Range("M134").Value = -62595849
$ws.Range("N134").Value = -15870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 6806074
$ws.Range("I132").Value = 8772648
$ws.Range("J132").Value = 12453.182
$ws.Range("K132").Value = 26317944
$ws.Range("L132").Value = 37359.546
$ws.Range("M132").Value = -26315414
$ws.Range("N132").Value = -42419.546

$ws.Range("H134").Value = 9471067
$ws.Range("I134").Value = 9260570
$ws.Range("J134").Value = 10418300
$ws.Range("K134").Value = 27781710
$ws.Range("L134").Value = 31254900
$ws.Range("M134").Value = -27779175
$ws.Range("N134").Value = -31259970

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 80913.16
$ws.Range("I137").Value = 103483.8
$ws.Range("J137").Value = 5677.6665
$ws.Range("K137").Value = 310451.4
$ws.Range("L137").Value = 17032.9995
$ws.Range("M137").Value = -305351.4
$ws.Range("N137").Value = -27232.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1431.6666
$ws.Range("I126").Value = 1458
$ws.Range("J126").Value = 1300
$ws.Range("K126").Value = 4374
$ws.Range("L126").Value = 3900
$ws.Range("M126").Value = -1904
$ws.Range("N126").Value = -8840

$ws.Range("H132").Value = 17564580
$ws.Range("I132").Value = 22751690
$ws.Range("J132").Value = 8210.308000000001
$ws.Range("K132").Value = 68255070
$ws.Range("L132").Value = 24630.924
$ws.Range("M132").Value = -68252540
$ws.Range("N132").Value = -29690.924

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3485.4055
$ws.Range("I40").Value = 3592.0967
$ws.Range("J40").Value = 2934.1667
$ws.Range("K40").Value = 3592.0967
$ws.Range("L40").Value = 2934.1667
$ws.Range("M40").Value = -3456.0967
$ws.Range("N40").Value = -3206.1667

$ws.Range("H122").Value = 10871478
$ws.Range("I122").Value = 1831.25
$ws.Range("K122").Value = 5493.75
$ws.Range("M122").Value = -3043.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1074
$ws.Range("I122").Value = 1029.5
$ws.Range("J122").Value = 1133.3334
$ws.Range("K122").Value = 3088.5
$ws.Range("L122").Value = 3400.0002
$ws.Range("M122").Value = -638.5
$ws.Range("N122").Value = -8300.0002

$ws.Range("H132").Value = 273721.84
$ws.Range("I132").Value = 24080.387
$ws.Range("J132").Value = 1272287.6
$ws.Range("K132").Value = 72241.16099999999
$ws.Range("L132").Value = 3816862.8
$ws.Range("M132").Value = -69711.16099999999
$ws.Range("N132").Value = -3821922.8
